$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 118 (the "このメスヒグマはベテランハンターだ" / bear entry),
# shifting all subsequent rows up by one.
$ws.Rows.Item(118).Delete()
